# Update 想去人数 (interest counts) in F column across sheets, per commit at 456a3b4
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 808
$ws.Range("F4").Value = 13524
$ws.Range("F5").Value = 13345
$ws.Range("F6").Value = 1032
$ws.Range("F8").Value = 24
$ws.Range("F9").Value = 578
$ws.Range("F11").Value = 11
$ws.Range("F12").Value = 16
$ws.Range("F13").Value = 712
$ws.Range("F14").Value = 2115
$ws.Range("F15").Value = 36
$ws.Range("F16").Value = 73
$ws.Range("F17").Value = 53
$ws.Range("F18").Value = 91
$ws.Range("F21").Value = 315
$ws.Range("F22").Value = 299
$ws.Range("F23").Value = 467
$ws.Range("F24").Value = 792
$ws.Range("F25").Value = 47
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 41
$ws.Range("F7").Value = 150
$ws.Range("F8").Value = 792
$ws.Range("F10").Value = 17
$ws.Range("F11").Value = 44
$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 77
$ws = $wb.Worksheets.Item(4)
$ws.Range("F5").Value = 808
$ws.Range("F6").Value = 13524
$ws.Range("F7").Value = 13345
$ws.Range("F8").Value = 1032
$ws.Range("F10").Value = 24
$ws.Range("F11").Value = 578
$ws.Range("F13").Value = 11
$ws.Range("F14").Value = 16
$ws.Range("F15").Value = 712
$ws.Range("F16").Value = 41
$ws.Range("F18").Value = 2115
$ws.Range("F19").Value = 36
$ws.Range("F20").Value = 73
$ws.Range("F21").Value = 53
$ws.Range("F22").Value = 91
$ws.Range("F26").Value = 77
$ws.Range("F28").Value = 315
$ws.Range("F29").Value = 299
$ws.Range("F30").Value = 467
$ws.Range("F31").Value = 792
$ws.Range("F32").Value = 150
$ws.Range("F33").Value = 792
$ws.Range("F35").Value = 17
$ws.Range("F36").Value = 47
$ws.Range("F37").Value = 44

$wb.Save()
